# Add columns I (I0) and J (IF) to the worksheet, mirroring the style of
# the existing header row (column H, style index 1: bold/centered/bordered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - same style as the other header cells (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I and J, rows 2-25
$values = @(
    @(6, 7),
    @(6, 8),
    @(5, 8),
    @(4, 8),
    @(5, 7),
    @(5, 7),
    @(5, 8),
    @(1, 3),
    @(10, 10),
    @(6, 8),
    @(5, 9),
    @(8, 9),
    @(6, 6),
    @(7, 7),
    @(9, 9),
    @(8, 9),
    @(3, 6),
    @(6, 8),
    @(3, 6),
    @(6, 8),
    @(7, 8),
    @(7, 9),
    @(3, 4),
    @(1, 2)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
